$wb = $excel.ActiveWorkbook

# --- "Delete" sheet: replace the old "Terminate after KA Module" value with
# "Order Take" in B2/D2 (this also makes the old shared string unused, so it
# gets dropped from the shared-string table on save, shifting later indices
# down by one - matching the rest of the workbook automatically). Also flip
# the On/Off flag in E2 and drop the top/bottom borders on that cell so it
# matches the "Off" look used elsewhere (style used by F2).
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("B2").Value = "Order Take"
$wsDelete.Range("D2").Value = "Order Take"
$wsDelete.Range("E2").Value = "Off"
$wsDelete.Range("E2").Borders.Item(8).LineStyle = -4142
$wsDelete.Range("E2").Borders.Item(9).LineStyle = -4142

# --- Update the remembered cursor/selection position on every sheet to
# reflect where the user last clicked before saving.
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Activate()
$wsCreate.Range("C17").Select()

$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Activate()
$wsEdit.Range("D15").Select()

$wsDelete.Activate()
$wsDelete.Range("D14").Select()

$wsQueries = $wb.Worksheets.Item("Queries")
$wsQueries.Activate()
$wsQueries.Range("A12").Select()
